$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.661.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -0.75%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.587.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.68%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.25%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'207.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -2.06%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  -3.62%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.31%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'22.22"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -4.59%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -2.25%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -2.85%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -1.47%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.813.49"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -2.60%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.589.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -2.24%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'  -4.04%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  -4.78%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'27.649.02"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -0.92%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'63.43"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'219.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -3.82%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.0₃0696"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -3.21%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  -4.12%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  +0.38%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  -4.78%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  -2.78%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -3.84%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'153.63"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.09%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'6.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -1.60%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.32%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'15.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -2.09%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  -4.87%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -2.43%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.0468"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -2.74%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  -5.28%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.371.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -3.35%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  -5.51%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -4.95%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -3.36%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  -0.67%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  -1.19%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.536"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -3.16%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -3.31%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +0.30%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.971"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.45%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'64.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -2.71%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +2.05%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'5.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -4.27%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'1.724.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -2.63%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  -5.07%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'87.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -1.47%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.0₆0101"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -1.28%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  -4.54%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0495"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.52%  "
$ws.Range("E51").ClearFormats()
